$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.418.03"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.642.48"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'529.84"
$ws.Range("D6").Value = "'145.45"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.571"
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "'0.339"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D13").Value = "3.107.96"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "59.391.17"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'21.01"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "2.643.42"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'343.48"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "'4.47"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "'10.64"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").Value = "'6.40"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'65.72"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").Value = "'0.420"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'7.27"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").Value = "0.0₃0804"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "'6.44"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").Value = "'19.11"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "'150.08"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "'4.21"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").Value = "'0.888"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "'0.869"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'3.67"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.0976"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.603"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("D44").Value = "'270.95"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "'19.46"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "2.038.96"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "'4.80"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'18.94"
$ws.Range("E51").Value = "  -0.35%  "
